# Change the table style on the table that lives on slide 16 of the
# "C1--C2-and-C3-PowerPoint" deck.
#
# Before:  a:tableStyleId = {F926386E-C90B-4B20-85DA-3244907729E6}
# After:   a:tableStyleId = {A76637DE-5832-48F9-AFE1-BAE8CDED8809}
#
# This mirrors picking a different style from the Table Styles gallery
# (Table Tools > Design) for the table's graphic frame.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(16)

# Find the shape that actually hosts the table instead of hard-coding an
# index, so this keeps working even if shape ordering ever shifts.
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{A76637DE-5832-48F9-AFE1-BAE8CDED8809}")
        break
    }
}
